$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the specific D/E cells we are about to overwrite
# so Excel keeps numeric-looking strings (e.g. "1.00", "532.38") as text,
# matching the source workbook where every value is stored as an inline string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.298.53"
$ws.Range("E2").Value = "  -3.52%  "

$ws.Range("D3").Value = "3.046.30"
$ws.Range("E3").Value = "  -3.07%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "532.38"
$ws.Range("E5").Value = "  -5.70%  "

$ws.Range("D6").Value = "131.17"
$ws.Range("E6").Value = "  -9.84%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.041.37"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("E11").Value = "  -10.23%  "

$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -4.10%  "

$ws.Range("D14").Value = "33.69"
$ws.Range("E14").Value = "  -8.87%  "

$ws.Range("D15").Value = "3.494.03"
$ws.Range("E15").Value = "  -4.26%  "

$ws.Range("D16").Value = "62.223.08"
$ws.Range("E16").Value = "  -3.73%  "

$ws.Range("E17").Value = "  -2.38%  "

$ws.Range("D18").Value = "3.041.67"
$ws.Range("E18").Value = "  -3.16%  "

$ws.Range("E19").Value = "  -5.55%  "

$ws.Range("D20").Value = "473.53"
$ws.Range("E20").Value = "  -8.03%  "

$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -6.80%  "

$ws.Range("D22").Value = "0.690"
$ws.Range("E22").Value = "  -4.06%  "

$ws.Range("E23").Value = "  -5.71%  "

$ws.Range("D24").Value = "78.20"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Value = "11.79"
$ws.Range("E25").Value = "  -8.65%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("E27").Value = "  -6.93%  "

$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  -10.26%  "

$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").Value = "25.47"
$ws.Range("E30").Value = "  -4.41%  "

$ws.Range("D31").Value = "1.83"
$ws.Range("E31").Value = "  -15.72%  "

$ws.Range("E32").Value = "  -4.72%  "

$ws.Range("E33").Value = "  -9.69%  "

$ws.Range("D34").Value = "56.30"
$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "5.82"

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.15"
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").Value = "465.09"
$ws.Range("E37").Value = "  -15.60%  "

$ws.Range("D38").Value = "3.068.78"
$ws.Range("E38").Value = "  -2.79%  "

$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").Value = "  -11.51%  "

$ws.Range("D40").Value = "0.0774"
$ws.Range("E40").Value = "  -6.35%  "

$ws.Range("D41").Value = "7.94"
$ws.Range("E41").Value = "  -4.32%  "

$ws.Range("E42").Value = "  -9.61%  "

$ws.Range("D43").Value = "2.53"
$ws.Range("E43").Value = "  -7.47%  "

$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "0.246"
$ws.Range("E45").Value = "  -8.26%  "

$ws.Range("E46").Value = "  -10.78%  "

$ws.Range("D47").Value = "0.0₃0511"
$ws.Range("E47").Value = "  -2.21%  "

$ws.Range("E48").Value = "  -6.76%  "

$ws.Range("E49").Value = "  -2.27%  "

$ws.Range("D50").Value = "115.16"
$ws.Range("E50").Value = "  -4.73%  "

$ws.Range("E51").Value = "  -8.10%  "

